$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 11).Value = 5508
$ws.Cells.Item(3, 11).Value = 5651
$ws.Cells.Item(4, 11).Value = 1180
$ws.Cells.Item(6, 11).Value = 6291
$ws.Cells.Item(7, 11).Value = 19034

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 11).Value = 349
$ws.Cells.Item(3, 11).Value = 384
$ws.Cells.Item(6, 11).Value = 428
$ws.Cells.Item(7, 11).Value = 1269

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(6, 11).Value = 98
$ws.Cells.Item(7, 11).Value = 424

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 11).Value = 222
$ws.Cells.Item(3, 11).Value = 300
$ws.Cells.Item(4, 11).Value = 39
$ws.Cells.Item(6, 11).Value = 238
$ws.Cells.Item(7, 11).Value = 817

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 11).Value = 108
$ws.Cells.Item(7, 11).Value = 322

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 11).Value = 209
$ws.Cells.Item(4, 11).Value = 31
$ws.Cells.Item(7, 11).Value = 641

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(6, 11).Value = 159
$ws.Cells.Item(7, 11).Value = 435

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 11).Value = 81
$ws.Cells.Item(7, 11).Value = 320

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(4, 11).Value = 72
$ws.Cells.Item(5, 11).Value = 46
$ws.Cells.Item(7, 11).Value = 559
$ws.Cells.Item(8, 11).Value = 1269
$ws.Cells.Item(10, 11).Value = 105
$ws.Cells.Item(15, 11).Value = 193
$ws.Cells.Item(17, 11).Value = 37
$ws.Cells.Item(18, 11).Value = 127
$ws.Cells.Item(19, 11).Value = 555
$ws.Cells.Item(20, 11).Value = 442
$ws.Cells.Item(22, 11).Value = 52
$ws.Cells.Item(23, 11).Value = 196
$ws.Cells.Item(24, 11).Value = 58
$ws.Cells.Item(29, 11).Value = 1026
$ws.Cells.Item(31, 11).Value = 204
$ws.Cells.Item(33, 11).Value = 817
$ws.Cells.Item(37, 11).Value = 641
$ws.Cells.Item(42, 11).Value = 707
$ws.Cells.Item(43, 11).Value = 166
$ws.Cells.Item(44, 11).Value = 165
$ws.Cells.Item(47, 11).Value = 132
$ws.Cells.Item(48, 11).Value = 242
$ws.Cells.Item(51, 11).Value = 239
$ws.Cells.Item(52, 11).Value = 495
$ws.Cells.Item(54, 11).Value = 367
$ws.Cells.Item(55, 11).Value = 211
$ws.Cells.Item(63, 11).Value = 56
$ws.Cells.Item(64, 11).Value = 123
$ws.Cells.Item(65, 11).Value = 435
$ws.Cells.Item(66, 11).Value = 63
$ws.Cells.Item(67, 11).Value = 723
$ws.Cells.Item(71, 11).Value = 59
$ws.Cells.Item(72, 11).Value = 91
$ws.Cells.Item(73, 11).Value = 167
$ws.Cells.Item(79, 11).Value = 483
$ws.Cells.Item(80, 11).Value = 67
$ws.Cells.Item(83, 11).Value = 424
$ws.Cells.Item(84, 11).Value = 147
$ws.Cells.Item(85, 11).Value = 896
$ws.Cells.Item(89, 11).Value = 277
$ws.Cells.Item(90, 11).Value = 173
$ws.Cells.Item(94, 11).Value = 258
$ws.Cells.Item(95, 11).Value = 322
$ws.Cells.Item(96, 11).Value = 204
$ws.Cells.Item(97, 11).Value = 153
$ws.Cells.Item(99, 11).Value = 320
$ws.Cells.Item(100, 11).Value = 37
$ws.Cells.Item(101, 11).Value = 19034

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(6, 11).Value = 73
$ws.Cells.Item(7, 11).Value = 204

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(6, 11).Value = 201
$ws.Cells.Item(7, 11).Value = 723

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(2, 11).Value = 51
$ws.Cells.Item(7, 11).Value = 147

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(3, 11).Value = 92
$ws.Cells.Item(7, 11).Value = 367

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 11).Value = 292
$ws.Cells.Item(3, 11).Value = 368
$ws.Cells.Item(7, 11).Value = 1026

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 11).Value = 33
$ws.Cells.Item(6, 11).Value = 119
$ws.Cells.Item(7, 11).Value = 242

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 11).Value = 166
$ws.Cells.Item(4, 11).Value = 24
$ws.Cells.Item(7, 11).Value = 555

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(3, 11).Value = 44
$ws.Cells.Item(7, 11).Value = 165

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(6, 11).Value = 266
$ws.Cells.Item(7, 11).Value = 707

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(3, 11).Value = 19
$ws.Cells.Item(7, 11).Value = 105

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(6, 11).Value = 73
$ws.Cells.Item(7, 11).Value = 211

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(2, 11).Value = 21
$ws.Cells.Item(7, 11).Value = 58

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(2, 11).Value = 54
$ws.Cells.Item(7, 11).Value = 196

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(3, 11).Value = 40
$ws.Cells.Item(7, 11).Value = 204

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 11).Value = 161
$ws.Cells.Item(6, 11).Value = 122
$ws.Cells.Item(7, 11).Value = 483

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(6, 11).Value = 46
$ws.Cells.Item(7, 11).Value = 123

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(3, 11).Value = 142
$ws.Cells.Item(6, 11).Value = 128
$ws.Cells.Item(7, 11).Value = 442

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(6, 11).Value = 35
$ws.Cells.Item(7, 11).Value = 127

$ws = $wb.Worksheets.Item('Burnside')
$ws.Cells.Item(6, 11).Value = 10
$ws.Cells.Item(7, 11).Value = 37

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Cells.Item(2, 11).Value = 8
$ws.Cells.Item(7, 11).Value = 37

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(3, 11).Value = 180
$ws.Cells.Item(6, 11).Value = 150
$ws.Cells.Item(7, 11).Value = 559

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(6, 11).Value = 113
$ws.Cells.Item(7, 11).Value = 258

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(3, 11).Value = 36
$ws.Cells.Item(7, 11).Value = 132

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 11).Value = 69
$ws.Cells.Item(7, 11).Value = 193

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(6, 11).Value = 31
$ws.Cells.Item(7, 11).Value = 63

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 11).Value = 55
$ws.Cells.Item(7, 11).Value = 167

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(6, 11).Value = 89
$ws.Cells.Item(7, 11).Value = 153

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(3, 11).Value = 86
$ws.Cells.Item(7, 11).Value = 277

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(6, 11).Value = 21
$ws.Cells.Item(7, 11).Value = 46

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(4, 11).Value = 14
$ws.Cells.Item(7, 11).Value = 173

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(2, 11).Value = 68
$ws.Cells.Item(3, 11).Value = 63
$ws.Cells.Item(7, 11).Value = 239

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(6, 11).Value = 66
$ws.Cells.Item(7, 11).Value = 166

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(3, 11).Value = 305
$ws.Cells.Item(6, 11).Value = 222
$ws.Cells.Item(7, 11).Value = 896

$ws = $wb.Worksheets.Item('Clearing')
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(7, 11).Value = 52

$ws = $wb.Worksheets.Item('Oakland')
$ws.Cells.Item(6, 11).Value = 15
$ws.Cells.Item(7, 11).Value = 59

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(6, 11).Value = 48
$ws.Cells.Item(7, 11).Value = 91

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Cells.Item(6, 11).Value = 32
$ws.Cells.Item(7, 11).Value = 67

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 11).Value = 133
$ws.Cells.Item(7, 11).Value = 495

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(6, 11).Value = 28
$ws.Cells.Item(7, 11).Value = 72
